$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 4-7 with new values (rows shift/rotate as report grows) ---

# Row 4
$ws.Range("B4").Value = 3
$ws.Range("C4").Value = "d5ed6f59-5f4c-40ef-a9aa-0d8f4c930b48"
$ws.Range("D4").Value = "BGOjoywll3sSUYo2pBwY"
$ws.Range("F4").Value = "Wm0PNk31XidAdvgB9iUo"

# Row 5
$ws.Range("B5").Value = 4
$ws.Range("D5").Value = "5dwojHXzFd2RzJ2ZPNtD"
$ws.Range("F5").Value = "YShT4GJbXiVYmuGPYfv4"

# Row 6
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "b129d86a-4077-4b49-8494-4942b76bfbb1"
$ws.Range("D6").Value = "OyhrkDZsWy64SkqfZAY5"
$ws.Range("E6").Value = 5
$ws.Range("F6").Value = "eurKXa2qR74n9sOpw4X0"

# Row 7
$ws.Range("B7").Value = 3
$ws.Range("D7").Value = "PqpflZDmB5tBiM15v9MQ"
$ws.Range("E7").Value = 4
$ws.Range("F7").Value = "hJMVI1CTbDSS4108H6Qp"

# --- Add new rows 8-12 ---

# Apply the same style used on column A (bold/bordered) to the new A cells
$ws.Range("A7").Copy()
$ws.Range("A8:A12").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 5
$ws.Range("C8").Value = "ae868375-47d2-4b8a-995b-66e8e33ebf0a"
$ws.Range("D8").Value = "BGOjoywll3sSUYo2pBwY"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = "hsCX7T7tqPKrlJp6WCcH"

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = "4825af3c-d463-4c5c-ac54-781f3890c126"
$ws.Range("D9").Value = "5dwojHXzFd2RzJ2ZPNtD"
$ws.Range("E9").Value = 4
$ws.Range("F9").Value = "iAQ9g8XGZ9ik4QY9CUdZ"

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 4
$ws.Range("C10").Value = "4825af3c-d463-4c5c-ac54-781f3890c126"
$ws.Range("D10").Value = "OyhrkDZsWy64SkqfZAY5"
$ws.Range("E10").Value = 4
$ws.Range("F10").Value = "tvllY9jEMZ4FEdMYjpzt"

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = 3
$ws.Range("C11").Value = "d5ed6f59-5f4c-40ef-a9aa-0d8f4c930b48"
$ws.Range("D11").Value = "PqpflZDmB5tBiM15v9MQ"
$ws.Range("E11").Value = 5
$ws.Range("F11").Value = "xRRIvdDlH0BG453SGVSz"

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = 4
$ws.Range("C12").Value = "ae868375-47d2-4b8a-995b-66e8e33ebf0a"
$ws.Range("D12").Value = "OyhrkDZsWy64SkqfZAY5"
$ws.Range("E12").Value = 2
$ws.Range("F12").Value = "yHra9FHoHYK939FWtJEc"
